$d = $word.ActiveDocument

$replacements = @(
    @("27÷9=", "51÷2="),
    @("49÷4=", "54÷7="),
    @("82÷7=", "89÷9="),
    @("83÷3=", "59÷5="),
    @("33÷6=", "65÷7="),
    @("19÷6=", "10÷7="),
    @("83÷7=", "95÷4="),
    @("43÷3=", "83÷6="),
    @("19÷5=", "59÷3="),
    @("44÷9=", "67÷4="),
    @("92÷2=", "51÷2="),
    @("47÷5=", "50÷4="),
    @("64÷2=", "51÷2="),
    @("64÷3=", "29÷8="),
    @("96÷4=", "59÷4="),
    @("20÷9=", "30÷2="),
    @("56÷8=", "98÷7="),
    @("14÷7=", "57÷5="),
    @("12÷7=", "20÷3="),
    @("70÷3=", "87÷9="),
    @("49÷5=", "83÷9="),
    @("37÷9=", "62÷4="),
    @("37÷8=", "57÷2="),
    @("36÷4=", "42÷3="),
    @("94÷2=", "45÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
